$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 419, shifting the existing rows 419-429 down to 420-430.
$ws.Rows("419:419").Insert()

# Populate the newly inserted row 419 with the new data record.
$ws.Range("A419").Value = 5
$ws.Range("B419").Value = "Macroferia Regional de Talca"
$ws.Range("C419").Value = "Maule"
$ws.Range("D419").Value = 45075
$ws.Range("E419").Value = 7
$ws.Range("F419").Value = 100112008
$ws.Range("G419").Value = "Coliflor"
$ws.Range("H419").Value = "Sin especificar"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 5000
$ws.Range("K419").Value = 500
$ws.Range("L419").Value = 500
$ws.Range("M419").Value = 500
$ws.Range("N419").Value = "`$/unidad"
$ws.Range("O419").Value = "Región del Maule"
$ws.Range("P419").Value = 500
$ws.Range("Q419").Value = 1
$ws.Range("R419").Value = "Hortaliza"
